$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.350.38"
$ws.Range("E2").Value = "  +0.74%  "

$ws.Range("D3").Value = "3.572.77"
$ws.Range("E3").Value = "  -0.08%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.24"
$ws.Range("E5").Value = "  +2.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.86"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").Value = "3.563.99"
$ws.Range("E7").Value = "  -0.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.620"
$ws.Range("E8").Value = "  -0.67%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.200"
$ws.Range("E10").Value = "  +8.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.648"
$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.78"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("E13").Value = "  +0.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.53"
$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("D15").Value = "4.140.20"
$ws.Range("E15").Value = "  -0.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.45"
$ws.Range("E16").Value = "  -1.08%  "

$ws.Range("D17").Value = "70.307.85"
$ws.Range("E17").Value = "  +0.55%  "

$ws.Range("D18").Value = "3.553.37"
$ws.Range("E18").Value = "  -1.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.53"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("E20").Value = "  -0.87%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "542.02"
$ws.Range("E21").Value = "  +10.79%  "

$ws.Range("E22").Value = "  -0.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.97"
$ws.Range("E23").Value = "  -7.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.71"
$ws.Range("E24").Value = "  +9.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.92"
$ws.Range("E25").Value = "  +0.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "95.83"
$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.48"
$ws.Range("E27").Value = "  +4.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.01"
$ws.Range("E28").Value = "  +2.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.15"
$ws.Range("E29").Value = "  -1.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.29"
$ws.Range("E30").Value = "  +1.42%  "

$ws.Range("E31").Value = "  -2.52%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.57"
$ws.Range("E32").Value = "  +4.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "65.26"
$ws.Range("E33").Value = "  -2.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.115"
$ws.Range("E34").Value = "  -0.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "553.90"
$ws.Range("E35").Value = "  -3.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.25"
$ws.Range("E36").Value = "  +6.58%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.417"
$ws.Range("E37").Value = "  +6.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.52"
$ws.Range("E38").Value = "  +1.12%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("E40").Value = "  -3.58%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.135"
$ws.Range("E41").Value = "  -1.53%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.371.57"
$ws.Range("E42").Value = "  +4.58%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("E43").Value = "  -4.53%  "

$ws.Range("E44").Value = "  -6.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.56"
$ws.Range("E45").Value = "  +3.58%  "

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0446"
$ws.Range("E47").Value = "  +2.20%  "

$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.22"
$ws.Range("E48").Value = "  -4.60%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.136"
$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("E50").Value = "  +0.02%  "

$ws.Range("B51").Value = "OceanProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.43"
$ws.Range("E51").Value = "  +18.41%  "
